$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 231
$ws.Range("F3").Value = 566
$ws.Range("F7").Value = 3139
$ws.Range("F8").Value = 2727
$ws.Range("F10").Value = 42
$ws.Range("F12").Value = 340
$ws.Range("F13").Value = 281
$ws.Range("F15").Value = 5610
$ws.Range("F17").Value = 1017
$ws.Range("F18").Value = 52
$ws.Range("F19").Value = 79
$ws.Range("F20").Value = 446
$ws.Range("F21").Value = 1207
$ws.Range("F23").Value = 108
$ws.Range("F24").Value = 327
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 239
$ws.Range("F9").Value = 50
$ws.Range("F18").Value = 51
$ws.Range("F24").Value = 289
$ws.Range("F33").Value = 40
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2550
$ws.Range("F6").Value = 1107
$ws.Range("F9").Value = 1425
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2550
$ws.Range("F6").Value = 1107
$ws.Range("F7").Value = 1425
$ws.Range("F11").Value = 231
$ws.Range("F12").Value = 566
$ws.Range("F15").Value = 3139
$ws.Range("F16").Value = 2727
$ws.Range("F18").Value = 42
$ws.Range("F20").Value = 239
$ws.Range("F22").Value = 340
$ws.Range("F24").Value = 50
$ws.Range("F26").Value = 5610
$ws.Range("F30").Value = 1017
$ws.Range("F32").Value = 52
$ws.Range("F33").Value = 79
$ws.Range("F36").Value = 51
$ws.Range("F40").Value = 1207
$ws.Range("F46").Value = 40
$ws.Range("F47").Value = 327
